# Apply updated cryptocurrency price (column D) and 1h volume/change (column E)
# values to Sheet1, as scraped by the GitHub Actions job that produced this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value. A leading "'" forces Excel to keep a
# numeric-looking Price string (e.g. "312.80") stored as text, matching the
# workbook's existing inline-string cells instead of letting Excel coerce it
# into a Double.
$updates = [ordered]@{
    "D2" = '27.843.65'
    "E2" = '  -0.80%  '
    "D3" = '1.902.32'
    "E3" = '  -0.34%  '
    "E4" = '  -0.42%  '
    "D5" = '''312.80'
    "E5" = '  -1.13%  '
    "D6" = '''1.002'
    "E6" = '  -0.46%  '
    "D7" = '''0.5011'
    "E7" = '  +3.75%  '
    "D8" = '''0.3811'
    "E8" = '  -0.22%  '
    "D9" = '''0.07292'
    "E9" = '  -0.96%  '
    "D10" = '''0.9091'
    "E10" = '  -2.60%  '
    "D11" = '''20.85'
    "E11" = '  +0.45%  '
    "D12" = '''0.07654'
    "E12" = '  -2.42%  '
    "D13" = '1.900.69'
    "E13" = '  -0.29%  '
    "D14" = '''5.478'
    "E14" = '  -0.41%  '
    "D15" = '''6.605'
    "E15" = '  -0.04%  '
    "D16" = '''91.28'
    "E16" = '  +0.09%  '
    "E17" = '  -0.44%  '
    "D18" = '''0.000008700'
    "E18" = '  -1.39%  '
    "E19" = '  -0.39%  '
    "D20" = '27.871.04'
    "E20" = '  -0.78%  '
    "D21" = '''14.51'
    "E21" = '  -2.03%  '
    "D22" = '''5.146'
    "E22" = '  -0.17%  '
    "D23" = '''10.82'
    "E23" = '  -0.35%  '
    "D24" = '''154.40'
    "E24" = '  -1.42%  '
    "D25" = '''1.859'
    "E25" = '  -3.33%  '
    "D26" = '''2.220'
    "E26" = '  +5.70%  '
    "E27" = '  -0.92%  '
    "D28" = '''115.13'
    "E28" = '  -1.04%  '
    "D29" = '''4.920'
    "E29" = '  -0.77%  '
    "D30" = '''0.08967'
    "E30" = '  +0.72%  '
    "D31" = '''3.219'
    "E31" = '  -4.23%  '
    "D32" = '''1.237'
    "E32" = '  -0.55%  '
    "D33" = '''0.7687'
    "E33" = '  +0.26%  '
    "D34" = '''4.634'
    "E34" = '  -0.85%  '
    "E35" = '  +0.73%  '
    "D36" = '''2.553'
    "E36" = '  -1.94%  '
    "D37" = '''1.098'
    "E37" = '  +0.11%  '
    "D38" = '''0.5534'
    "E38" = '  +0.94%  '
    "E39" = '  +1.01%  '
    "D40" = '''0.05267'
    "D41" = '''6.962'
    "E41" = '  -0.69%  '
    "D42" = '''8.531'
    "E42" = '  +1.02%  '
    "D43" = '''0.1522'
    "E43" = '  -0.02%  '
    "D44" = '''111.21'
    "E44" = '  +3.95%  '
    "D45" = '''10.61'
    "E45" = '  -0.92%  '
    "D46" = '''0.4791'
    "E46" = '  -0.78%  '
    "D47" = '''1.002'
    "E47" = '  -0.55%  '
    "D48" = '''1.635'
    "E48" = '  -1.20%  '
    "D49" = '''67.33'
    "E49" = '  -1.50%  '
    "E50" = '  -0.29%  '
    "D51" = '''0.9004'
    "E51" = '  -0.35%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
